$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 4289.313007116892
$ws.Range("C3").Value = 4280.169215615671
$ws.Range("C4").Value = 4280.169215615671
$ws.Range("C5").Value = 4280.169215615671
$ws.Range("C6").Value = 4090.507436366653
$ws.Range("C7").Value = 4090.507436366653
$ws.Range("C8").Value = 4090.507436366653
$ws.Range("C9").Value = 4076.13586166109
$ws.Range("C10").Value = 4076.13586166109
$ws.Range("C11").Value = 4076.13586166109
$ws.Range("C12").Value = 4076.13586166109
